$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 9 with FALCON_F2 strategy, copying the layout/style pattern
# used by the existing FALCON_F1 row (row 8).
$ws.Range("A9").Value = 39
$ws.Range("B9").Value = "FALCON_F2"
$ws.Range("C9").Value = "Trade Entry and Treshold is predicted from AI"
$ws.Range("D9").Value = "read prediction from AI and enter if it is matching the market type"
$ws.Range("E9").Value = "when prediction time has expired or targets were reached"
$ws.Range("F9").Value = "Undefined. Using Pattern recognition NN"

# Match formatting (wrap text + row height) of the row above it
$ws.Range("C9:F9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 45

# Update selection to reflect the newly entered row, matching the
# post-edit cursor position left by the author.
$ws.Range("D9:F9").Select()
